$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Standard..." report-name strings (added to the shared string table first)
$ws.Range("G2").Value = "StandardExcelReport-Takeda - MM Maintenance-Clinical-2023_"
$ws.Range("G5").Value = "StandardExcelReport-Takeda - MM Maintenance-Economic-2023_"
$ws.Range("G8").Value = "StandardExcelReport-Takeda - MM Maintenance-Quality of Life-2023_"
$ws.Range("G11").Value = "StandardExcelReport-Takeda - MM Maintenance-Real-world Evidence-2023_"

# New "ExcelReport-Takeda-..." strings (added to the shared string table next)
$ws.Range("G3").Value = "ExcelReport-Takeda-MM Maintenance-Clinical-"
$ws.Range("G6").Value = "ExcelReport-Takeda-MM Maintenance-Economic-"
$ws.Range("G9").Value = "ExcelReport-Takeda-MM Maintenance-Quality of Life-"
$ws.Range("G12").Value = "ExcelReport-Takeda-MM Maintenance-Real-world Evidence-"

# Reuse of existing "WordReport-Takeda - MM Maintenance-..." strings, shifted up a row each
$ws.Range("G4").Value = "WordReport-Takeda - MM Maintenance-Clinical-"
$ws.Range("G7").Value = "WordReport-Takeda - MM Maintenance-Economic-"
$ws.Range("G10").Value = "WordReport-Takeda - MM Maintenance-Quality of Life-"
$ws.Range("G13").Value = "WordReport-Takeda - MM Maintenance-Real-world Evidence-"

# Remove the now-unused trailing rows
$ws.Range("G14:G18").ClearContents() | Out-Null

$ws.Range("G9").Select() | Out-Null
